$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 203, shifting existing rows 203:224 down to 204:225
$ws.Rows("203:203").Insert()

# Populate the newly inserted row 203 with the new record
$ws.Range("A203").Value = 4
$ws.Range("B203").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C203").Value = "Los Lagos"
$ws.Range("D203").Value = 44449
$ws.Range("D203").NumberFormat = $ws.Range("D204").NumberFormat
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = 100114001
$ws.Range("G203").Value = "Papa"
$ws.Range("H203").Value = "Patagonia"
$ws.Range("I203").Value = "1a (guarda)"
$ws.Range("J203").Value = 600
$ws.Range("K203").Value = 7000
$ws.Range("L203").Value = 7500
$ws.Range("M203").Value = 7250
$ws.Range("N203").Value = "`$/saco 25 kilos"
$ws.Range("O203").Value = "Provincia de Llanquihue"
$ws.Range("P203").Value = 290
$ws.Range("Q203").Value = 25
$ws.Range("R203").Value = "Hortaliza"
